# Insert a new weekly data row at row 4 (shifts existing rows 4-26 down to 5-27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44503
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 72
$ws.Range("K4").Value = 1600
$ws.Range("L4").Value = 1600
$ws.Range("M4").Value = 1600
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 1600
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
